$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246, shifting existing rows 246-328 down to 247-329
$ws.Rows(246).Insert()

# Populate the newly inserted row 246 with the new weekly price record
$ws.Range("A246").Value = 3
$ws.Range("B246").Value = "Femacal de La Calera"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = 44588
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E246").Value = 5
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100108
$ws.Range("H246").Value = "Tropicales y subtropicales"
$ws.Range("I246").Value = 100108002
$ws.Range("J246").Value = "Mango"
$ws.Range("K246").Value = "Sin especificar"
$ws.Range("L246").Value = "Primera"
$ws.Range("M246").Value = 260
$ws.Range("N246").Value = 6500
$ws.Range("O246").Value = 7000
$ws.Range("P246").Value = 6731
$ws.Range("Q246").Value = "$/bandeja 4 kilos"
$ws.Range("R246").Value = "Perú"
$ws.Range("S246").Value = 1683
$ws.Range("T246").Value = 4
